$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registration")

# Update the "Regression" column (C) for rows 6 through 69 to "NO",
# matching the "Functional" column (D) values - removing the unwanted
# library from the Regression test register.
$ws.Range("C6:C69").Value = "NO"

# Reflect the selection/active cell change recorded in the saved file.
$ws.Activate()
$ws.Range("C6:C69").Select()
